$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a serial date value that was bumped
# from 45175 (2023-09-06) to 45177 (2023-09-08) for every data row
# (rows 2 through 359).
$ws.Range("C2:C359").Value = 45177
